$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").Value = "  "
$ws.Range("B22").Value = "Саитов Артур (Филипов...)"

$ws.Range("H26").Value = 5
$ws.Range("I26").Value = 5
$ws.Range("I26").Style = $ws.Range("H24").Style
$ws.Range("J26").Formula = "=SUM(C26:I26)"

$ws.Range("P11").Select()
